$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 2.471608527004719;  E = 3.410691294819117;  F = 0.2358042635023594;  G = 0.7053456474095583 }
    3  = @{ D = 2.338538888841867;  E = 2.43381213163957;   F = 0.1692694444209337;  G = 0.216906065819785 }
    4  = @{ D = 4.675471302121878;  E = 4.276767794974148;  F = 0.1688678255304694;  G = 0.06919194874353707 }
    5  = @{ D = 5.168602275429294;  E = 5.393329305574298;  F = 0.03372045508585871; G = 0.07866586111485958 }
    6  = @{ D = 3.855129688978195;  E = 3.243498809635639;  F = 0.2850432296593984;  G = 0.0811662698785464 }
    7  = @{ D = 2.411616187542677;  E = 2.285308588296175;  F = 0.2058080937713385;  G = 0.1426542941480875 }
    8  = @{ D = 4.397048187442124;  E = 4.716597873717546;  F = 0.09926204686053097; G = 0.1791494684293866 }
    9  = @{ D = 2.329123958013952;  E = 2.737986572086811;  F = 0.1645619790069759;  G = 0.3689932860434055 }
    10 = @{ D = 4.394832290709019;  E = 3.822696465998888;  F = 0.4649440969030063;  G = 0.2742321553329627 }
    11 = @{ D = 3.868999440222979;  E = 3.620934154838324;  F = 0.2896664800743262;  G = 0.2069780516127745 }
    12 = @{ D = 3.50879477057606;   E = 3.620779767632484;  F = 0.1695982568586866;  G = 0.2069265892108282 }
    13 = @{ D = 4.320871399249882;  E = 4.570783686824143;  F = 0.08021784981247038; G = 0.1426959217060357 }
    14 = @{ D = 3.01353463344276;   E = 2.814172249287367;  F = 0.5067673167213798;  G = 0.4070861246436834 }
    15 = @{ D = 3.618761103600264;  E = 3.926798224449158;  F = 0.2062537012000879;  G = 0.3089327414830526 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
